$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance ("1") in column L for every row that already has other
# weekly attendance marks, except rows 6, 13 and 14 (those stay blank in L).
$attendedRows = @(2,3,4,5,7,8,9,10,11,12,15,16,17,18,19,20)
foreach ($r in $attendedRows) {
    $ws.Cells.Item($r, 12).Value = 1
}

# New cell comments describing the reasons behind specific marks.
$c13 = $ws.Range("L13").AddComment()
[void]$c13.Text("수학축전")

$c14 = $ws.Range("L14").AddComment()
[void]$c14.Text("결막염 진단서 제출해야함")

$c15 = $ws.Range("L15").AddComment()
[void]$c15.Text("지각")

# Update the active selection to match the cell that was last interacted with.
[void]$ws.Range("N11").Select()
